$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C31").Value = "RoleDaoTest 空指针异常"
$ws.Range("D31").Value = "23:30--23:45"

$ws.Range("D31").Select()
